$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 480, pushing existing rows 480-599 down to 481-600
$ws.Rows("480:480").Insert()

# Populate the new row 480 with the new data record
$ws.Range("A480").Value = 10
$ws.Range("B480").Value = "Vega Modelo de Temuco"
$ws.Range("C480").Value = "La Araucanía"
$ws.Range("D480").Value = 44782
$ws.Range("E480").Value = 9
$ws.Range("F480").Value = 100112027
$ws.Range("G480").Value = "Melón"
$ws.Range("H480").Value = "Calameño"
$ws.Range("I480").Value = "Primera"
$ws.Range("J480").Value = 80
$ws.Range("K480").Value = 23000
$ws.Range("L480").Value = 23000
$ws.Range("M480").Value = 23000
$ws.Range("N480").Value = "$/caja 16 unidades"
$ws.Range("O480").Value = "Brasil"
$ws.Range("P480").Value = 1438
$ws.Range("Q480").Value = 16
$ws.Range("R480").Value = "Hortaliza"

# Match the date number format used in column D elsewhere
$ws.Range("D480").NumberFormat = $ws.Range("D481").NumberFormat()
